# The workbook stores, for each genome bin (rows 2-4, column A), a predicted
# placement ("1-o__Elusimicrobiales") together with a numeric score in column B
# and a supporting count in column C. This update refreshes the previously
# placeholder score values in column B with the recomputed scores from the
# updated "ful-path.csv" pipeline run, without altering any of the labels,
# headers, or other data in the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 1643.9452004067166
$ws.Range("B3").Value = 1535.0922864223048
$ws.Range("B4").Value = 1723.472007440319
